# Applies the Jan 25 2024 cryptos-list refresh (prices + 1h volume %,
# plus a couple of rank swaps) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '39.872.01'
$ws.Cells.Item(2, 5).Value = '  -0.64%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.207.31'
$ws.Cells.Item(3, 5).Value = '  -1.26%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.14%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '291.53'
$ws.Cells.Item(5, 5).Value = '  -0.29%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '86.52'
$ws.Cells.Item(6, 5).Value = '  -0.50%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.509'
$ws.Cells.Item(7, 5).Value = '  -1.41%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.04%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.467'
$ws.Cells.Item(9, 5).Value = '  -1.81%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '30.21'
$ws.Cells.Item(10, 5).Value = '  -3.25%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0778'
$ws.Cells.Item(11, 5).Value = '  -1.39%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'OKB'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '49.99'
$ws.Cells.Item(12, 5).Value = '  +5.77%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.111'
$ws.Cells.Item(13, 5).Value = '  +2.48%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.46'
$ws.Cells.Item(14, 5).Value = '  +1.32%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.543.57'
$ws.Cells.Item(15, 5).Value = '  -1.49%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '13.69'
$ws.Cells.Item(16, 5).Value = '  -3.44%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.227.34'
$ws.Cells.Item(17, 5).Value = '  -0.70%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.729'
$ws.Cells.Item(18, 5).Value = '  -0.41%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '39.766.48'
$ws.Cells.Item(19, 5).Value = '  -0.74%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -0.85%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '11.18'
$ws.Cells.Item(21, 5).Value = '  -0.79%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.73'
$ws.Cells.Item(22, 5).Value = '  -1.74%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '65.31'
$ws.Cells.Item(23, 5).Value = '  -0.86%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '237.87'
$ws.Cells.Item(24, 5).Value = '  +0.39%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -0.03%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -1.53%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '1.82'
$ws.Cells.Item(27, 5).Value = '  -2.00%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '23.43'
$ws.Cells.Item(28, 5).Value = '  +2.31%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Cosmos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '9.19'
$ws.Cells.Item(29, 5).Value = '  -1.25%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '2.05'
$ws.Cells.Item(30, 5).Value = '  -7.34%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '156.32'
$ws.Cells.Item(31, 5).Value = '  +2.88%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '31.41'
$ws.Cells.Item(32, 5).Value = '  -5.67%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.03%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.94'
$ws.Cells.Item(34, 5).Value = '  -0.84%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0707'
$ws.Cells.Item(35, 5).Value = '  -2.29%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.88'
$ws.Cells.Item(36, 5).Value = '  +2.23%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'WEMIXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.32'
$ws.Cells.Item(37, 5).Value = '  -2.41%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0976'
$ws.Cells.Item(39, 5).Value = '  -2.47%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '15.28'
$ws.Cells.Item(40, 5).Value = '  -5.70%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.69'
$ws.Cells.Item(41, 5).Value = '  -1.53%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '2.122.51'
$ws.Cells.Item(42, 5).Value = '  +2.57%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -2.34%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0268'
$ws.Cells.Item(44, 5).Value = '  -0.56%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -2.53%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '9.80'
$ws.Cells.Item(46, 5).Value = '  -1.44%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '17.49'
$ws.Cells.Item(47, 5).Value = '  -6.04%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.67'
$ws.Cells.Item(48, 5).Value = '  +1.95%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.415.49'
$ws.Cells.Item(49, 5).Value = '  -1.56%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.48'
$ws.Cells.Item(50, 5).Value = '  +2.17%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '88.48'
$ws.Cells.Item(51, 5).Value = '  -1.29%  '
